$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Align rows 2 and 3 with row 4's account/loan data (tarjetas no propias stabilization) ---

# Column D (usuario) - copy value+style from D4
$ws.Range("D4").Copy($ws.Range("D2"))
$ws.Range("D4").Copy($ws.Range("D3"))

# Column E (clave) - copy value+style from E4
$ws.Range("E4").Copy($ws.Range("E2"))
$ws.Range("E4").Copy($ws.Range("E3"))

# Column O (tipoPrestamo) - copy value+style from O4
$ws.Range("O4").Copy($ws.Range("O2"))
$ws.Range("O4").Copy($ws.Range("O3"))

# Column P (numeroPrestamo) - copy value+style from P4
$ws.Range("P4").Copy($ws.Range("P2"))
$ws.Range("P4").Copy($ws.Range("P3"))

# Column S (tipoCuenta) - row 3 becomes "Corriente" like the others
$ws.Range("S3").Value = $ws.Range("S2").Value()

# Column T (numeroCuenta) - copy value+style from T4 so all rows share the same account number/style
$ws.Range("T4").Copy($ws.Range("T2"))
$ws.Range("T4").Copy($ws.Range("T3"))

# --- Reset the view: scroll back to A1 and move the active selection to D4 ---
$ws.Range("D4").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
